{"js": "const replacements = [\n    [\"2025-01-20 Monday\", \"2025-01-21 Tuesday\"],\n    [\"41\u00d727=\", \"20\u00d791=\"],\n    [\"77\u00d719=\", \"34\u00d753=\"],\n    [\"40\u00d787=\", \"23\u00d771=\"],\n    [\"46\u00d722=\", \"42\u00d758=\"],\n    [\"96\u00d790=\", \"58\u00d745=\"],\n    [\"73\u00d764=\", \"43\u00d722=\"],\n    [\"89\u00d761=\", \"37\u00d745=\"],\n    [\"20\u00d793=\", \"13\u00d757=\"],\n    [\"96\u00d769=\", \"44\u00d763=\"],\n    [\"69\u00d790=\", \"38\u00d721=\"],\n    [\"24\u00d746=\", \"20\u00d767=\"],\n    [\"22\u00d777=\", \"49\u00d724=\"],\n    [\"80\u00d725=\", \"54\u00d750=\"],\n    [\"90\u00d763=\", \"75\u00d711=\"],\n    [\"89\u00d777=\", \"74\u00d734=\"],\n    [\"82\u00d760=\", \"50\u00d790=\"],\n    [\"76\u00d731=\", \"90\u00d717=\"],\n    [\"69\u00d764=\", \"61\u00d796=\"],\n    [\"92\u00d739=\", \"70\u00d721=\"],\n    [\"97\u00d797=\", \"91\u00d765=\"],\n    [\"43\u00d777=\", \"98\u00d781=\"],\n    [\"11\u00d799=\", \"31\u00d798=\"],\n    [\"36\u00d715=\", \"27\u00d754=\"],\n    [\"21\u00d725=\", \"22\u00d757=\"],\n    [\"39\u00d741=\", \"80\u00d760=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old=\"2025-01-20 Monday\"; new=\"2025-01-21 Tuesday\"},\n    @{old=\"41\u00d727=\"; new=\"20\u00d791=\"},\n    @{old=\"77\u00d719=\"; new=\"34\u00d753=\"},\n    @{old=\"40\u00d787=\"; new=\"23\u00d771=\"},\n    @{old=\"46\u00d722=\"; new=\"42\u00d758=\"},\n    @{old=\"96\u00d790=\"; new=\"58\u00d745=\"},\n    @{old=\"73\u00d764=\"; new=\"43\u00d722=\"},\n    @{old=\"89\u00d761=\"; new=\"37\u00d745=\"},\n    @{old=\"20\u00d793=\"; new=\"13\u00d757=\"},\n    @{old=\"96\u00d769=\"; new=\"44\u00d763=\"},\n    @{old=\"69\u00d790=\"; new=\"38\u00d721=\"},\n    @{old=\"24\u00d746=\"; new=\"20\u00d767=\"},\n    @{old=\"22\u00d777=\"; new=\"49\u00d724=\"},\n    @{old=\"80\u00d725=\"; new=\"54\u00d750=\"},\n    @{old=\"90\u00d763=\"; new=\"75\u00d711=\"},\n    @{old=\"89\u00d777=\"; new=\"74\u00d734=\"},\n    @{old=\"82\u00d760=\"; new=\"50\u00d790=\"},\n    @{old=\"76\u00d731=\"; new=\"90\u00d717=\"},\n    @{old=\"69\u00d764=\"; new=\"61\u00d796=\"},\n    @{old=\"92\u00d739=\"; new=\"70\u00d721=\"},\n    @{old=\"97\u00d797=\"; new=\"91\u00d765=\"},\n    @{old=\"43\u00d777=\"; new=\"98\u00d781=\"},\n    @{old=\"11\u00d799=\"; new=\"31\u00d798=\"},\n    @{old=\"36\u00d715=\"; new=\"27\u00d754=\"},\n    @{old=\"21\u00d725=\"; new=\"22\u00d757=\"},\n    @{old=\"39\u00d741=\"; new=\"80\u00d760=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.Text = $r.old\n    $find.Replacement.Text = $r.new\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n"}
